# Commit: Mon, Jul 27, 2020  1:05:25 PM
#
# The presentation's tables were re-styled: the three slides that contain a
# table (14, 15, 16) had their table style switched from the deck's custom
# "Table_0" style ({902330B2-8F6A-4A41-A4E1-59F472CBE8ED}) to the built-in
# "No Style, Table Grid" style ({AEEE09BA-4F7E-4508-BE53-50B1DCA18EA0}).

$p = $ppt.ActivePresentation

$newStyleId = "{AEEE09BA-4F7E-4508-BE53-50B1DCA18EA0}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
